$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Mark rows 37-40 ("A faire" column C) by copying the task label from
# column B into column C and highlighting it with a yellow fill plus a
# thick left border, matching the styling already used elsewhere in the
# sheet for completed/flagged status cells.
foreach ($r in 37..40) {
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value2 = $cellB.Value2
    $cellC.Interior.Color = 65535
    $cellC.Borders.Item(7).Weight = 4
}

$ws.Range("C40").Select()
